$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1 -- copy formatting (bold, border, alignment) from the
# neighboring header cell G1, then set the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New numeric value in H2
$ws.Range("H2").Value = 1
